$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Functions table (column D = CallbackId, column E = Description) ---
$ws.Range("D11").Value = 6
$ws.Range("D11").HorizontalAlignment = -4108
$ws.Range("E11").Value = "Mouse wheel moved"

# --- Events table (column G = EventId, column H = Description) ---
$ws.Range("G6").Value = 1
$ws.Range("G6").HorizontalAlignment = -4108
$ws.Range("H6").Value = "Zoom canvas"

$ws.Range("G7").Value = 2
$ws.Range("G7").HorizontalAlignment = -4108
$ws.Range("H7").Value = "Move vertex"

$ws.Range("G8").Value = 3
$ws.Range("G8").HorizontalAlignment = -4108
$ws.Range("H8").Value = "Created isolated node"

$ws.Range("G9").Value = 4
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("H9").Value = "shift-click vertex"

# --- Selection / active cell ---
$ws.Range("H15").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
